# Applies the three paragraph-formatting changes described by the diff:
#  1) Paragraph 1: right-align + add 567-twip right indent.
#  2) Paragraph 2: add 1.5 line spacing (line=360, lineRule=auto).
#  3) Paragraph 3: replace the "TEST TEST TEST" runs (and their proofErr
#     spell-check wrappers) with a single "Style 2" run in Times New Roman
#     11pt, and drop the bold paragraph-mark formatting.

$d = $word.ActiveDocument

# --- Paragraph 1: "This header is ok and bold is allowed" ---------------
$p1 = $d.Paragraphs(1)
$p1.Alignment = 2            # wdAlignParagraphRight
$p1.Format.RightIndent = 28.35   # points == 567 twips

# --- Paragraph 2: "This header does not have bold" ----------------------
$p2 = $d.Paragraphs(2)
$p2.Format.LineSpacingRule = 1   # wdLineSpace1pt5 -> line=360 lineRule=auto

# --- Paragraph 3: " TEST TEST TEST" -> "Style 2" -------------------------
$p3 = $d.Paragraphs(3)
$start3 = $p3.Range.Start
$end3 = $p3.Range.End

# Remove the old paragraph content (runs + proofErr marks) entirely.
$d.Range($start3, $end3).Delete()

# Re-insert a clean paragraph with the exact target formatting in its
# place, preserving the original paragraph identity attributes so the
# surrounding markup stays stable.
$lastP = $d.Paragraphs($d.Paragraphs.Count)
$splitAt = $lastP.Range.End

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$newParaXml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document ' + $xmlNs + '><w:body>' +
  '<w:p w14:paraId="3113E59B" w14:textId="42DFFF9F" w:rsidR="00B17C0B" w:rsidRPr="00B17C0B" w:rsidRDefault="00B17C0B" w:rsidP="002B65ED">' +
  '<w:pPr><w:jc w:val="center"/><w:rPr>' +
  '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
  '<w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>Style 2</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Range($splitAt - 1, $splitAt - 1).InsertXML($newParaXml)
